$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 50019296
$ws.Range("I62").Value = 100005850
$ws.Range("K62").Value = 100005850
$ws.Range("M62").Value = -100005226
$ws.Range("H65").Value = 50019296
$ws.Range("I65").Value = 100005850
$ws.Range("K65").Value = 500029250
$ws.Range("M65").Value = -500026130
$ws.Range("H98").Value = 33234.4
$ws.Range("I98").Value = 33511.35
$ws.Range("K98").Value = 33511.35
$ws.Range("M98").Value = -32013.35
$ws.Range("H122").Value = 33234.4
$ws.Range("I122").Value = 33511.35
$ws.Range("K122").Value = 100534.05
$ws.Range("M122").Value = -98084.04999999999
$ws.Range("H138").Value = 2812.68
$ws.Range("J138").Value = 4880.5415
$ws.Range("L138").Value = 14641.6245
$ws.Range("N138").Value = -24921.6245
$ws.Range("H140").Value = 82833
$ws.Range("J140").Value = 82833
$ws.Range("L140").Value = 82833
$ws.Range("N140").Value = -93193
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6732.5386
$ws.Range("I45").Value = 6504.1665
$ws.Range("K45").Value = 6504.1665
$ws.Range("M45").Value = -6127.1665
$ws.Range("H61").Value = 5370.6
$ws.Range("I61").Value = 4813.3335
$ws.Range("K61").Value = 4813.3335
$ws.Range("M61").Value = -4601.3335
$ws.Range("H132").Value = 4810.857
$ws.Range("I132").Value = 4156.737
$ws.Range("K132").Value = 12470.211
$ws.Range("M132").Value = -9940.210999999999
$ws.Range("H136").Value = 5370.6
$ws.Range("I136").Value = 4813.3335
$ws.Range("K136").Value = 14440.0005
$ws.Range("M136").Value = -11890.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2317.6562
$ws.Range("I134").Value = 1845.3928
$ws.Range("K134").Value = 5536.178400000001
$ws.Range("M134").Value = -3001.178400000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2969.4
$ws.Range("I31").Value = 1242
$ws.Range("K31").Value = 1242
$ws.Range("M31").Value = -947
$ws.Range("H34").Value = 2969.4
$ws.Range("I34").Value = 1242
$ws.Range("K34").Value = 1242
$ws.Range("M34").Value = -1040
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 2288.6667
$ws.Range("J94").Value = 1974
$ws.Range("L94").Value = 1974
$ws.Range("N94").Value = -2876
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3025.6667
$ws.Range("I3").Value = 2557.75
$ws.Range("J3").Value = 3400
$ws.Range("K3").Value = 7673.25
$ws.Range("L3").Value = 10200
$ws.Range("M3").Value = -7561.25
$ws.Range("N3").Value = -10424
$ws.Range("H34").Value = 1002042.5
$ws.Range("I34").Value = 2274749.5
$ws.Range("J34").Value = 2058.5
$ws.Range("K34").Value = 6824248.5
$ws.Range("L34").Value = 6175.5
$ws.Range("M34").Value = -6824164.5
$ws.Range("N34").Value = -6343.5
$ws.Range("H39").Value = 1057.5714
$ws.Range("I39").Value = 247.9
$ws.Range("J39").Value = 3081.75
$ws.Range("K39").Value = 743.7
$ws.Range("L39").Value = 9245.25
$ws.Range("M39").Value = -449.7
$ws.Range("N39").Value = -9833.25
$ws.Range("H55").Value = 9045.516
$ws.Range("I55").Value = 923.25
$ws.Range("J55").Value = 10165.827
$ws.Range("K55").Value = 2769.75
$ws.Range("L55").Value = 30497.481
$ws.Range("M55").Value = -2592.75
$ws.Range("N55").Value = -30851.481
$ws.Range("H86").Value = 1206.091
$ws.Range("I86").Value = 1159.4
$ws.Range("K86").Value = 3478.2
$ws.Range("M86").Value = -2292.2
$ws.Range("H89").Value = 1206.091
$ws.Range("I89").Value = 1159.4
$ws.Range("K89").Value = 10434.6
$ws.Range("M89").Value = -4506.6
$ws.Range("H139").Value = 1502473.1
$ws.Range("I139").Value = 1877216.5
$ws.Range("K139").Value = 5631649.5
$ws.Range("M139").Value = -5626509.5
$ws.Range("H140").Value = 11388.765
$ws.Range("I140").Value = 11850.5625
$ws.Range("K140").Value = 35551.6875
$ws.Range("M140").Value = -30371.6875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 9000
$ws.Range("J21").Value = 9000
$ws.Range("L21").Value = 9000
$ws.Range("N21").Value = -9346
$ws.Range("H30").Value = 9000
$ws.Range("J30").Value = 9000
$ws.Range("L30").Value = 9000
$ws.Range("N30").Value = -9210
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 71000
$ws.Range("J125").Value = 71000
$ws.Range("L125").Value = 71000
$ws.Range("N125").Value = -80840
$ws.Range("H136").Value = 8074.3213
$ws.Range("I136").Value = 14749
$ws.Range("K136").Value = 44247
$ws.Range("M136").Value = -41697
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 39995
$ws.Range("J21").Value = 39995
$ws.Range("L21").Value = 39995
$ws.Range("N21").Value = -40465
$ws.Range("H26").Value = 20011.5
$ws.Range("I26").Value = 20011.5
$ws.Range("K26").Value = 20011.5
$ws.Range("M26").Value = -19718.5
$ws.Range("H35").Value = 39995
$ws.Range("J35").Value = 39995
$ws.Range("L35").Value = 39995
$ws.Range("N35").Value = -40575
$ws.Range("H75").Value = 12500
$ws.Range("I75").Value = 12500
$ws.Range("K75").Value = 12500
$ws.Range("M75").Value = -11564
$ws.Range("H78").Value = 12500
$ws.Range("I78").Value = 12500
$ws.Range("K78").Value = 37500
$ws.Range("M78").Value = -32820
$ws.Range("H81").Value = 25028.428
$ws.Range("I81").Value = 37800
$ws.Range("J81").Value = 7999.6665
$ws.Range("K81").Value = 75600
$ws.Range("L81").Value = 15999.333
$ws.Range("M81").Value = -74539
$ws.Range("N81").Value = -18121.333
$ws.Range("H84").Value = 25028.428
$ws.Range("I84").Value = 37800
$ws.Range("J84").Value = 7999.6665
$ws.Range("K84").Value = 378000
$ws.Range("L84").Value = 79996.66500000001
$ws.Range("M84").Value = -372696
$ws.Range("N84").Value = -90604.66500000001
$ws.Range("H107").Value = 32566.8
$ws.Range("I107").Value = 3547.8462
$ws.Range("J107").Value = 86459.14
$ws.Range("K107").Value = 10643.5386
$ws.Range("L107").Value = 259377.42
$ws.Range("M107").Value = -8723.5386
$ws.Range("N107").Value = -263217.42

Write-Host "Applied market-price refresh across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
